$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column G ("Hora") changes from 18 to 19 for every data row (2-51).
# Values are prefixed with an apostrophe so Excel stores them as text
# (matching the workbook's existing text storage for this column) instead
# of silently coercing the numeric-looking string to a number.
$ws.Range("G2:G51").Value = "'19"

# Row 2
$ws.Range("D2").Value = "'236.76"

# Row 3
$ws.Range("D3").Value = "'21.88"

# Row 4
$ws.Range("B4").Value = 'LEO'
$ws.Range("C4").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D4").Value = "'3.918"
$ws.Range("E4").Value = '3LEOLEO'

# Row 5
$ws.Range("B5").Value = 'HuobiToken'
$ws.Range("C5").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D5").Value = "'5.347"
$ws.Range("E5").Value = '4HuobiTokenHT'

# Row 6
$ws.Range("B6").Value = 'Cronos'
$ws.Range("C6").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D6").Value = "'0.05584"
$ws.Range("E6").Value = '5CronosCRO'

# Row 7
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").Value = "'6.456"
$ws.Range("E7").Value = '6KuCoinTokenKCS'

# Row 8
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").Value = "'3.357"
$ws.Range("E8").Value = '7GateTokenGT'

# Row 9
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = "'0.7993"
$ws.Range("E9").Value = '8MXTokenMX'

# Row 10
$ws.Range("B10").Value = 'FTXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D10").Value = "'1.041"
$ws.Range("E10").Value = '9FTXTokenFTT'

# Row 11
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = "'0.1396"
$ws.Range("E11").Value = '10WazirXWRX'

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = "'0.07252"
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'

# Row 13
$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").Value = "'0.03192"
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'

# Row 14
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = "'0.02933"
$ws.Range("E14").Value = '13BitrueCoinBTR'

# Row 15
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = "'0.09235"
$ws.Range("E15").Value = '14BitMartTokenBMX'

# Row 16
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = "'0.001668"
$ws.Range("E16").Value = '15BitForexTokenBF'

# Row 17
$ws.Range("B17").Value = 'MCDex'
$ws.Range("C17").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D17").Value = "'3.255"
$ws.Range("E17").Value = '16MCDexMCB'

# Row 18
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").Value = "'0.04775"
$ws.Range("E18").Value = '17CoinExTokenCET'

# Row 19
$ws.Range("D19").Value = "'0.006259"

# Row 20
$ws.Range("D20").Value = "'0.005082"

# Row 24
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").Value = "'2.199"
$ws.Range("E24").Value = '23BTSETokenBTSE'

# Row 25
$ws.Range("B25").Value = 'One'
$ws.Range("C25").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D25").Value = "'0.01167"
$ws.Range("E25").Value = '24OneONEBestin24h'

# Row 27
$ws.Range("D27").Value = "'0.1256"

# Row 40
$ws.Range("D40").Value = "'0.04108"

# Row 41
$ws.Range("D41").Value = "'0.007074"

# Row 42
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = "'0.003499"
$ws.Range("E42").Value = '41CEJICEJI'

# Row 43
$ws.Range("B43").Value = 'BKEXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D43").Value = "'0.1037"
$ws.Range("E43").Value = '42BKEXTokenBKK'

# Row 44
$ws.Range("D44").Value = "'0.008940"

# Row 45
$ws.Range("D45").Value = "'0.00005440"

# Row 48
$ws.Range("D48").Value = "'0.03381"

# Row 49
$ws.Range("D49").Value = "'0.00002100"
